$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row
$ws.Range("A1").Value = 'Cluster Name'
$ws.Range("B1").Value = 'Active Cases'

# Data rows (2..64), sorted per final layout, with updated case counts
$ws.Range("A2").Value = '3323 Villa Maria Catholic Homes St Bernadette''s Aged Care Sunshine North'
$ws.Range("B2").Value = 15
$ws.Range("A3").Value = '3376 Royal Freemasons Coppin Centre Melbourne'
$ws.Range("B3").Value = 32
$ws.Range("A4").Value = '3601 Baptcare Westhaven community'
$ws.Range("B4").Value = 13
$ws.Range("A5").Value = '3653 Fronditha Thalpori St Albans Aged Care'
$ws.Range("B5").Value = 35
$ws.Range("A6").Value = '3825 TLC Forest Lodge Residential Aged Care Frankston North'
$ws.Range("B6").Value = 16
$ws.Range("A7").Value = '4167 Royal Freemasons Centennial Lodge Wantirna South'
$ws.Range("B7").Value = 18
$ws.Range("A8").Value = '44054 Dromana Primary School Dromana'
$ws.Range("B8").Value = 10
$ws.Range("A9").Value = '44121 Wallan Primary School Wallan'
$ws.Range("B9").Value = 10
$ws.Range("A10").Value = '44226 Boneo Primary School Boneo'
$ws.Range("B10").Value = 12
$ws.Range("A11").Value = '44234 Lucknow Primary School'
$ws.Range("B11").Value = 13
$ws.Range("A12").Value = '44321 Maiden Gully Primary School Maiden Gully'
$ws.Range("B12").Value = 11
$ws.Range("A13").Value = '44395 Buln Buln Primary School'
$ws.Range("B13").Value = 10
$ws.Range("A14").Value = '44852 Dandenong South Primary School Dandenong'
$ws.Range("B14").Value = 16
$ws.Range("A15").Value = '44978 Deer Park West Primary School Deer Park'
$ws.Range("B15").Value = 10
$ws.Range("A16").Value = '45034 River Gum Primary School Hampton Park'
$ws.Range("B16").Value = 11
$ws.Range("A17").Value = '45158 Rowellyn Park Primary School Carrum Downs'
$ws.Range("B17").Value = 10
$ws.Range("A18").Value = '45573 Narre Warren South P-12 College Narre Warren South'
$ws.Range("B18").Value = 16
$ws.Range("A19").Value = '45585 Mount Ridley College Craigieburn'
$ws.Range("B19").Value = 16
$ws.Range("A20").Value = '45695 Sacred Heart Primary School Yarrawonga'
$ws.Range("B20").Value = 53
$ws.Range("A21").Value = '4574 Village Glen Aged Care Residences Mornington'
$ws.Range("B21").Value = 13
$ws.Range("A22").Value = '45804 St Therese''s School Essendon'
$ws.Range("B22").Value = 14
$ws.Range("A23").Value = '45809 St Finbar''s Primary School Brighton East'
$ws.Range("B23").Value = 11
$ws.Range("A24").Value = '46050 Our Lady''s Catholic Primary School Craigieburn'
$ws.Range("B24").Value = 22
$ws.Range("A25").Value = '46322 Minaret College Officer Campus Officer'
$ws.Range("B25").Value = 40
$ws.Range("A26").Value = '46328 Ilim College Dallas Primary Campus Tier 1A Dallas'
$ws.Range("B26").Value = 10
$ws.Range("A27").Value = '46390 Al Siraat College Epping'
$ws.Range("B27").Value = 42
$ws.Range("A28").Value = '50176 Meadows Primary School Broadmeadows'
$ws.Range("B28").Value = 13
$ws.Range("A29").Value = '50516 Ilim College Glenroy Campus Hadfield'
$ws.Range("B29").Value = 13
$ws.Range("A30").Value = '50567 Alamanda K9 College Point Cook'
$ws.Range("B30").Value = 11
$ws.Range("A31").Value = '51478 Wyndham Vale Primary School Wyndham'
$ws.Range("B31").Value = 11
$ws.Range("A32").Value = '52380 Al Iman College Melton South'
$ws.Range("B32").Value = 31
$ws.Range("A33").Value = '52786 Hume Anglican Grammar Donnybrook Campus'
$ws.Range("B33").Value = 11
$ws.Range("A34").Value = '52912 Edgars Creek Primary School Wollert'
$ws.Range("B34").Value = 13
$ws.Range("A35").Value = '52985 Minaret College Springvale'
$ws.Range("B35").Value = 15
$ws.Range("A36").Value = 'Adass Israel School Elsternwick'
$ws.Range("B36").Value = 19
$ws.Range("A37").Value = 'Covenant College Bell Post Hill'
$ws.Range("B37").Value = 17
$ws.Range("A38").Value = 'Creekside K-9 College Caroline Springs'
$ws.Range("B38").Value = 18
$ws.Range("A39").Value = 'Derrimut Primary School Derrimut'
$ws.Range("B39").Value = 11
$ws.Range("A40").Value = 'Exford Primary School Exford'
$ws.Range("B40").Value = 15
$ws.Range("A41").Value = 'Hazel Glen College Doreen'
$ws.Range("B41").Value = 15
$ws.Range("A42").Value = 'Hazelwood North Primary School Hazelwood North'
$ws.Range("B42").Value = 11
$ws.Range("A43").Value = 'Ilim College Dallas Main Campus Dallas Oct'
$ws.Range("B43").Value = 22
$ws.Range("A44").Value = 'Ilim College Kiewa Campus Dallas'
$ws.Range("B44").Value = 11
$ws.Range("A45").Value = 'InverlochKongwak Primary School'
$ws.Range("B45").Value = 10
$ws.Range("A46").Value = 'Islamic College of Melbourne Tarneit Oct Nov'
$ws.Range("B46").Value = 53
$ws.Range("A47").Value = 'Lavalla Catholic College St Pauls Campus Traralgon'
$ws.Range("B47").Value = 13
$ws.Range("A48").Value = 'Middle Park Primary School Middle Park'
$ws.Range("B48").Value = 14
$ws.Range("A49").Value = 'Morwell Park Primary School Morwell Outbreak'
$ws.Range("B49").Value = 66
$ws.Range("A50").Value = 'Nio Early Learning Adventures Preston'
$ws.Range("B50").Value = 21
$ws.Range("A51").Value = 'Pentland Primary School Darley'
$ws.Range("B51").Value = 13
$ws.Range("A52").Value = 'Sirius College Shepparton Campus Shepparton'
$ws.Range("B52").Value = 14
$ws.Range("A53").Value = 'Society Restaurant Melbourne'
$ws.Range("B53").Value = 40
$ws.Range("A54").Value = 'St Ambrose Parish Primary School Woodend'
$ws.Range("B54").Value = 12
$ws.Range("A55").Value = 'St Georges Road Primary School Shepparton'
$ws.Range("B55").Value = 14
$ws.Range("A56").Value = 'St Paul''s Primary School Sunshine West'
$ws.Range("B56").Value = 11
$ws.Range("A57").Value = 'Stockdale Road Primary School Traralgon'
$ws.Range("B57").Value = 22
$ws.Range("A58").Value = 'Supreme Caravans Manufacturing Campbellfield'
$ws.Range("B58").Value = 52
$ws.Range("A59").Value = 'Top Yard Rooftop Melbourne'
$ws.Range("B59").Value = 13
$ws.Range("A60").Value = 'Truganina P-9 College Truganina'
$ws.Range("B60").Value = 15
$ws.Range("A61").Value = 'Tucker Road Bentleigh Primary School Bentleigh'
$ws.Range("B61").Value = 10
$ws.Range("A62").Value = 'Warragul Regional College Warragul'
$ws.Range("B62").Value = 20
$ws.Range("A63").Value = 'Werribee Mercy Hospital Emergency Department'
$ws.Range("B63").Value = 10
$ws.Range("A64").Value = 'Yeshivah College St Kilda East'
$ws.Range("B64").Value = 20
